$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("uncertainty_avoidance", 98, 0.096, 0.002, 11.486, 3.522, 3.262, 0.047, 25.703, 12.759, 2.014),
    @("indulgence", 79, 0.041, 0.069, 9.223, 4.995, 1.847, 0.57, 10.485, 18.364, 0.571),
    @("individualism", 98, 0.024, 0.119, 6.236, 3.966, 1.572, 0.226, 17.499, 14.368, 1.218),
    @("power_distance", 98, 0.024, 0.121, -5.742, 3.674, -1.563, 0, 87.39, 13.313, 6.564),
    @("democracy", 139, 0.011, 0.209, 0.417, 0.331, 1.261, 0.002, 3.854, 1.194, 3.229),
    @("masculinity", 55, 0.028, 0.22, -5.332, 4.3, -1.24, 0, 69.316, 15.374, 4.509),
    @("unemployment_rate", 62, 0.023, 0.237, 1.346, 1.126, 1.195, 0.703, 1.615, 4.212, 0.383),
    @("freedom_index", 163, 0.007, 0.272, 4.164, 3.781, 1.101, 0.001, 44.43, 13.471, 3.298),
    @("expected_years_of_schooling", 164, 0.004, 0.381, 0.317, 0.361, 0.878, 0, 12.37, 1.283, 9.638),
    @("motivation", 98, 0.004, 0.523, -1.996, 3.113, -0.641, 0, 54.222, 11.279, 4.808),
    @("mean_years_of_schooling", 163, 0.002, 0.588, -0.217, 0.4, -0.542, 0, 9.678, 1.423, 6.799),
    @("gdp_per_capita", 165, 0.001, 0.643, 0.097, 0.209, 0.464, 0, 8.448, 0.746, 11.319),
    @("corruption_score", 153, 0.001, 0.648, -1.09, 2.386, -0.457, 0, 46.904, 8.549, 5.487),
    @("safety_index", 116, 0.001, 0.678, -1.105, 2.651, -0.417, 0, 58.798, 9.528, 6.171),
    @("crime_index", 116, 0.001, 0.678, 1.105, 2.651, 0.417, 0, 41.202, 9.528, 4.324),
    @("long_term_orientation", 89, 0.001, 0.805, 0.995, 4.018, 0.248, 0.018, 35.299, 14.685, 2.404),
    @("human_development_index", 163, 0, 0.985, 0, 0.019, 0.019, 0, 0.717, 0.067, 10.651),
    @("life_expectancy", 165, 0, 0.99, -0.012, 0.972, -0.012, 0, 71.167, 3.458, 20.582),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($rowNum, 1).Value = $rowData[0]
    for ($j = 1; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = [double]$rowData[$j]
    }
}

$wb.Save()